$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the second data row (row 2) from the sheet; all rows below
# shift up by one, matching the target layout.
$ws.Rows.Item(2).EntireRow.Delete()
